$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 5333.3335
$ws.Range("J48").Value = 5333.3335
$ws.Range("L48").Value = 16000.0005
$ws.Range("N48").Value = -16584.0005
$ws.Range("H56").Value = 5333.3335
$ws.Range("J56").Value = 5333.3335
$ws.Range("L56").Value = 16000.0005
$ws.Range("N56").Value = -17068.0005
$ws.Range("H70").Value = 1999.3334
$ws.Range("H73").Value = 1999.3334
$ws.Range("H97").Value = 645.13336
$ws.Range("J97").Value = 655.5
$ws.Range("L97").Value = 1966.5
$ws.Range("N97").Value = -2958.5
$ws.Range("H138").Value = 3037.6042
$ws.Range("I138").Value = 2938.9
$ws.Range("J138").Value = 3063.5789
$ws.Range("K138").Value = 8816.700000000001
$ws.Range("L138").Value = 9190.736699999999
$ws.Range("M138").Value = -3676.700000000001
$ws.Range("N138").Value = -19470.7367
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15937.73
$ws.Range("I32").Value = 9814.983
$ws.Range("J32").Value = 29276.572
$ws.Range("K32").Value = 9814.983
$ws.Range("L32").Value = 29276.572
$ws.Range("M32").Value = -9527.983
$ws.Range("N32").Value = -29850.572
$ws.Range("H45").Value = 1401.3043
$ws.Range("I45").Value = 1301.7646
$ws.Range("J45").Value = 1683.3334
$ws.Range("K45").Value = 1301.7646
$ws.Range("L45").Value = 1683.3334
$ws.Range("M45").Value = -924.7646
$ws.Range("N45").Value = -2437.3334
$ws.Range("H61").Value = 333334660
$ws.Range("I61").Value = 333334660
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 333334660
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -333334448
$ws.Range("N61").ClearContents()
$ws.Range("H68").Value = 30099
$ws.Range("J68").Value = 30099
$ws.Range("L68").Value = 30099
$ws.Range("N68").Value = -31721
$ws.Range("H71").Value = 30099
$ws.Range("J71").Value = 30099
$ws.Range("L71").Value = 90297
$ws.Range("N71").Value = -98409
$ws.Range("H86").Value = 20000
$ws.Range("J86").Value = 20000
$ws.Range("L86").Value = 20000
$ws.Range("N86").Value = -22372
$ws.Range("H89").Value = 20000
$ws.Range("J89").Value = 20000
$ws.Range("L89").Value = 60000
$ws.Range("N89").Value = -71856
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 333334660
$ws.Range("I136").Value = 333334660
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 1000003980
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1000001430
$ws.Range("N136").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6939.4443
$ws.Range("I134").Value = 1084
$ws.Range("J134").Value = 11623.8
$ws.Range("K134").Value = 3252
$ws.Range("L134").Value = 34871.39999999999
$ws.Range("M134").Value = -717
$ws.Range("N134").Value = -39941.39999999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1637.4286
$ws.Range("I31").Value = 1605.4814
$ws.Range("J31").Value = 2500
$ws.Range("K31").Value = 1605.4814
$ws.Range("L31").Value = 2500
$ws.Range("M31").Value = -1310.4814
$ws.Range("N31").Value = -3090
$ws.Range("H34").Value = 1637.4286
$ws.Range("I34").Value = 1605.4814
$ws.Range("J34").Value = 2500
$ws.Range("K34").Value = 1605.4814
$ws.Range("L34").Value = 2500
$ws.Range("M34").Value = -1403.4814
$ws.Range("N34").Value = -2904
$ws.Range("H52").Value = 33647
$ws.Range("J52").Value = 35803.332
$ws.Range("L52").Value = 35803.332
$ws.Range("N52").Value = -36391.332
$ws.Range("H132").Value = 1944.4857
$ws.Range("I132").Value = 1505.3704
$ws.Range("K132").Value = 4516.1112
$ws.Range("M132").Value = -1986.1112
$ws.Range("H138").Value = 117597.14
$ws.Range("J138").Value = 117597.14
$ws.Range("L138").Value = 117597.14
$ws.Range("N138").Value = -127877.14
$ws.Range("H141").Value = 285263.3
$ws.Range("J141").Value = 285263.3
$ws.Range("L141").Value = 285263.3
$ws.Range("N141").Value = -295623.3
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4283931.5
$ws.Range("I4").Value = 1666704
$ws.Range("J4").Value = 6527269.5
$ws.Range("K4").Value = 5000112
$ws.Range("L4").Value = 19581808.5
$ws.Range("M4").Value = -5000000
$ws.Range("N4").Value = -19582032.5
$ws.Range("H9").Value = 1891
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 1891
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 5673
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -6121
$ws.Range("H70").Value = 5350.125
$ws.Range("I70").Value = 3444
$ws.Range("J70").Value = 5622.4287
$ws.Range("K70").Value = 10332
$ws.Range("L70").Value = 16867.2861
$ws.Range("M70").Value = -10017
$ws.Range("N70").Value = -17497.2861
$ws.Range("H73").Value = 5350.125
$ws.Range("I73").Value = 3444
$ws.Range("J73").Value = 5622.4287
$ws.Range("K73").Value = 10332
$ws.Range("L73").Value = 16867.2861
$ws.Range("M73").Value = -9240
$ws.Range("N73").Value = -19051.2861
$ws.Range("H74").Value = 3404.5
$ws.Range("J74").Value = 3956.25
$ws.Range("L74").Value = 11868.75
$ws.Range("N74").Value = -13990.75
$ws.Range("H77").Value = 3404.5
$ws.Range("J77").Value = 3956.25
$ws.Range("L77").Value = 35606.25
$ws.Range("N77").Value = -46214.25
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("N87").ClearContents()
$ws.Range("H88").Value = 5430
$ws.Range("J88").Value = 5760.7144
$ws.Range("L88").Value = 17282.1432
$ws.Range("N88").Value = -18138.1432
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("N90").ClearContents()
$ws.Range("H91").Value = 5430
$ws.Range("J91").Value = 5760.7144
$ws.Range("L91").Value = 17282.1432
$ws.Range("N91").Value = -20246.1432
$ws.Range("H131").Value = 21309416
$ws.Range("J131").Value = 38489.574
$ws.Range("L131").Value = 115468.722
$ws.Range("N131").Value = -125548.722
$ws.Range("H140").Value = 26439.428
$ws.Range("I140").Value = 52371.9
$ws.Range("J140").Value = 2864.4546
$ws.Range("K140").Value = 157115.7
$ws.Range("L140").Value = 8593.363799999999
$ws.Range("M140").Value = -151935.7
$ws.Range("N140").Value = -18953.3638
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 31145.428
$ws.Range("J134").Value = 31145.428
$ws.Range("L134").Value = 93436.284
$ws.Range("N134").Value = -98506.284
$ws.Range("H136").Value = 27365.1
$ws.Range("J136").Value = 27365.1
$ws.Range("L136").Value = 82095.29999999999
$ws.Range("N136").Value = -87195.29999999999
$ws.Range("H138").Value = 38680
$ws.Range("J138").Value = 38680
$ws.Range("L138").Value = 38680
$ws.Range("N138").Value = -48960
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4209.222
$ws.Range("J132").Value = 3713.7144
$ws.Range("L132").Value = 11141.1432
$ws.Range("N132").Value = -16201.1432
$ws.Range("H139").Value = 32727.5
$ws.Range("J139").Value = 32727.5
$ws.Range("L139").Value = 32727.5
$ws.Range("N139").Value = -43007.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 1266.6666
$ws.Range("J2").Value = 1266.6666
$ws.Range("L2").Value = 1266.6666
$ws.Range("N2").Value = -1490.6666
$ws.Range("H46").Value = 47874.5
$ws.Range("J46").Value = 47874.5
$ws.Range("L46").Value = 47874.5
$ws.Range("N46").Value = -48336.5
$ws.Range("H98").Value = 15000
$ws.Range("J98").Value = 15000
$ws.Range("L98").Value = 15000
$ws.Range("N98").Value = -20990
$ws.Range("H107").Value = 409
$ws.Range("I107").Value = 342.4
$ws.Range("J107").Value = 520
$ws.Range("K107").Value = 1027.2
$ws.Range("L107").Value = 1560
$ws.Range("M107").Value = 892.8000000000002
$ws.Range("N107").Value = -5400
$ws.Range("H134").Value = 47874.5
$ws.Range("J134").Value = 47874.5
$ws.Range("L134").Value = 143623.5
$ws.Range("N134").Value = -148693.5
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
